$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap data between paired rows (B..AC), column A (id) stays fixed per row ---
# Row 49
$ws.Cells.Item(49,2).Value = 6865311
$ws.Cells.Item(49,3).Value = "Bosnia Herzegovina Premier Liga"
$ws.Cells.Item(49,4).Value = "Bosnia  Herzegovina Premier Liga"
$ws.Cells.Item(49,5).Value = 45200.41666666666
$ws.Cells.Item(49,6).Value = "Sloga"
$ws.Cells.Item(49,7).Value = "GOSK Gabela"
$ws.Cells.Item(49,8).Value = 3
$ws.Cells.Item(49,9).Value = 2
$ws.Cells.Item(49,10).Value = "H"
$ws.Cells.Item(49,11).Value = 1.833
$ws.Cells.Item(49,12).Value = 3.6
$ws.Cells.Item(49,13).Value = 3.4
$ws.Cells.Item(49,14).Value = 1.909
$ws.Cells.Item(49,15).Value = 3.4
$ws.Cells.Item(49,16).Value = 3.3
$ws.Cells.Item(49,17).Value = -0.5
$ws.Cells.Item(49,18).Value = 1.925
$ws.Cells.Item(49,19).Value = 1.875
$ws.Cells.Item(49,20).Value = 2.25
$ws.Cells.Item(49,21).Value = 1.825
$ws.Cells.Item(49,22).Value = 1.975
$ws.Cells.Item(49,23).Value = 0.909
$ws.Cells.Item(49,24).Value = -1
$ws.Cells.Item(49,25).Value = -1
$ws.Cells.Item(49,26).Value = 0.925
$ws.Cells.Item(49,27).Value = -1
$ws.Cells.Item(49,28).Value = 0.825
$ws.Cells.Item(49,29).Value = -1

# Row 50
$ws.Cells.Item(50,2).Value = 6865310
$ws.Cells.Item(50,3).Value = "Bosnia Herzegovina Premier Liga"
$ws.Cells.Item(50,4).Value = "Bosnia  Herzegovina Premier Liga"
$ws.Cells.Item(50,5).Value = 45200.41666666666
$ws.Cells.Item(50,6).Value = "NK Igman Konjic"
$ws.Cells.Item(50,7).Value = "Zrinjski Mostar"
$ws.Cells.Item(50,8).Value = 0
$ws.Cells.Item(50,9).Value = 2
$ws.Cells.Item(50,10).Value = "A"
$ws.Cells.Item(50,11).Value = 3.4
$ws.Cells.Item(50,12).Value = 3.6
$ws.Cells.Item(50,13).Value = 1.833
$ws.Cells.Item(50,14).Value = 4.75
$ws.Cells.Item(50,15).Value = 4.75
$ws.Cells.Item(50,16).Value = 1.45
$ws.Cells.Item(50,17).Value = 1.25
$ws.Cells.Item(50,18).Value = 1.775
$ws.Cells.Item(50,19).Value = 2.025
$ws.Cells.Item(50,20).Value = 2.75
$ws.Cells.Item(50,21).Value = 1.85
$ws.Cells.Item(50,22).Value = 1.95
$ws.Cells.Item(50,23).Value = -1
$ws.Cells.Item(50,24).Value = -1
$ws.Cells.Item(50,25).Value = 0.45
$ws.Cells.Item(50,26).Value = -1
$ws.Cells.Item(50,27).Value = 1.025
$ws.Cells.Item(50,28).Value = -1
$ws.Cells.Item(50,29).Value = 0.95

# Row 76
$ws.Cells.Item(76,2).Value = 6865328
$ws.Cells.Item(76,3).Value = "Bosnia Herzegovina Premier Liga"
$ws.Cells.Item(76,4).Value = "Bosnia  Herzegovina Premier Liga"
$ws.Cells.Item(76,5).Value = 45235.375
$ws.Cells.Item(76,6).Value = "Siroki Brijeg"
$ws.Cells.Item(76,7).Value = "NK Posusje"
$ws.Cells.Item(76,8).Value = 1
$ws.Cells.Item(76,9).Value = 1
$ws.Cells.Item(76,10).Value = "D"
$ws.Cells.Item(76,11).Value = 2
$ws.Cells.Item(76,12).Value = 3
$ws.Cells.Item(76,13).Value = 3.5
$ws.Cells.Item(76,14).Value = 2.1
$ws.Cells.Item(76,15).Value = 3
$ws.Cells.Item(76,16).Value = 3.3
$ws.Cells.Item(76,17).Value = -0.25
$ws.Cells.Item(76,18).Value = 1.825
$ws.Cells.Item(76,19).Value = 1.975
$ws.Cells.Item(76,20).Value = 2
$ws.Cells.Item(76,21).Value = 1.825
$ws.Cells.Item(76,22).Value = 1.975
$ws.Cells.Item(76,23).Value = -1
$ws.Cells.Item(76,24).Value = 2
$ws.Cells.Item(76,25).Value = -1
$ws.Cells.Item(76,26).Value = -0.5
$ws.Cells.Item(76,27).Value = 0.4875
$ws.Cells.Item(76,28).Value = 0
$ws.Cells.Item(76,29).Value = -0

# Row 77
$ws.Cells.Item(77,2).Value = 6865377
$ws.Cells.Item(77,3).Value = "Bosnia Herzegovina Premier Liga"
$ws.Cells.Item(77,4).Value = "Bosnia  Herzegovina Premier Liga"
$ws.Cells.Item(77,5).Value = 45235.375
$ws.Cells.Item(77,6).Value = "Zrinjski Mostar"
$ws.Cells.Item(77,7).Value = "FK Tuzla City"
$ws.Cells.Item(77,8).Value = 3
$ws.Cells.Item(77,9).Value = 1
$ws.Cells.Item(77,10).Value = "H"
$ws.Cells.Item(77,11).Value = 1.333
$ws.Cells.Item(77,12).Value = 5
$ws.Cells.Item(77,13).Value = 6
$ws.Cells.Item(77,14).Value = 1.166
$ws.Cells.Item(77,15).Value = 6.5
$ws.Cells.Item(77,16).Value = 13
$ws.Cells.Item(77,17).Value = -2
$ws.Cells.Item(77,18).Value = 1.9
$ws.Cells.Item(77,19).Value = 1.9
$ws.Cells.Item(77,20).Value = 3.25
$ws.Cells.Item(77,21).Value = 1.95
$ws.Cells.Item(77,22).Value = 1.85
$ws.Cells.Item(77,23).Value = 0.1659999999999999
$ws.Cells.Item(77,24).Value = -1
$ws.Cells.Item(77,25).Value = -1
$ws.Cells.Item(77,26).Value = 0
$ws.Cells.Item(77,27).Value = -0
$ws.Cells.Item(77,28).Value = 0.95
$ws.Cells.Item(77,29).Value = -1

# Row 87
$ws.Cells.Item(87,2).Value = 7505497
$ws.Cells.Item(87,3).Value = "Bosnia Herzegovina Premier Liga"
$ws.Cells.Item(87,4).Value = "Bosnia  Herzegovina Premier Liga"
$ws.Cells.Item(87,5).Value = 45256.375
$ws.Cells.Item(87,6).Value = "Zeljeznicar"
$ws.Cells.Item(87,7).Value = "NK Posusje"
$ws.Cells.Item(87,8).Value = 1
$ws.Cells.Item(87,9).Value = 1
$ws.Cells.Item(87,10).Value = "D"
$ws.Cells.Item(87,11).Value = 1.65
$ws.Cells.Item(87,12).Value = 3.4
$ws.Cells.Item(87,13).Value = 4.75
$ws.Cells.Item(87,14).Value = 1.8
$ws.Cells.Item(87,15).Value = 3.2
$ws.Cells.Item(87,16).Value = 4.2
$ws.Cells.Item(87,17).Value = -0.5
$ws.Cells.Item(87,18).Value = 1.825
$ws.Cells.Item(87,19).Value = 1.975
$ws.Cells.Item(87,20).Value = 2
$ws.Cells.Item(87,21).Value = 1.75
$ws.Cells.Item(87,22).Value = 2.05
$ws.Cells.Item(87,23).Value = -1
$ws.Cells.Item(87,24).Value = 2.2
$ws.Cells.Item(87,25).Value = -1
$ws.Cells.Item(87,26).Value = -1
$ws.Cells.Item(87,27).Value = 0.9750000000000001
$ws.Cells.Item(87,28).Value = 0
$ws.Cells.Item(87,29).Value = -0

# Row 88
$ws.Cells.Item(88,2).Value = 7505495
$ws.Cells.Item(88,3).Value = "Bosnia Herzegovina Premier Liga"
$ws.Cells.Item(88,4).Value = "Bosnia  Herzegovina Premier Liga"
$ws.Cells.Item(88,5).Value = 45256.375
$ws.Cells.Item(88,6).Value = "Sloga"
$ws.Cells.Item(88,7).Value = "Zvijezda 09"
$ws.Cells.Item(88,8).Value = 1
$ws.Cells.Item(88,9).Value = 0
$ws.Cells.Item(88,10).Value = "H"
$ws.Cells.Item(88,11).Value = 1.444
$ws.Cells.Item(88,12).Value = 4.2
$ws.Cells.Item(88,13).Value = 5.5
$ws.Cells.Item(88,14).Value = 1.5
$ws.Cells.Item(88,15).Value = 4.2
$ws.Cells.Item(88,16).Value = 5.25
$ws.Cells.Item(88,17).Value = -1
$ws.Cells.Item(88,18).Value = 1.8
$ws.Cells.Item(88,19).Value = 2
$ws.Cells.Item(88,20).Value = 2.75
$ws.Cells.Item(88,21).Value = 1.775
$ws.Cells.Item(88,22).Value = 2.025
$ws.Cells.Item(88,23).Value = 0.5
$ws.Cells.Item(88,24).Value = -1
$ws.Cells.Item(88,25).Value = -1
$ws.Cells.Item(88,26).Value = 0
$ws.Cells.Item(88,27).Value = -0
$ws.Cells.Item(88,28).Value = -1
$ws.Cells.Item(88,29).Value = 1.025

# Row 111
$ws.Cells.Item(111,2).Value = 6865352
$ws.Cells.Item(111,3).Value = "Bosnia Herzegovina Premier Liga"
$ws.Cells.Item(111,4).Value = "Bosnia  Herzegovina Premier Liga"
$ws.Cells.Item(111,5).Value = 45339.375
$ws.Cells.Item(111,6).Value = "NK Posusje"
$ws.Cells.Item(111,7).Value = "Zvijezda 09"
$ws.Cells.Item(111,8).Value = 2
$ws.Cells.Item(111,9).Value = 0
$ws.Cells.Item(111,10).Value = "H"
$ws.Cells.Item(111,11).Value = 1.615
$ws.Cells.Item(111,12).Value = 3.5
$ws.Cells.Item(111,13).Value = 4.75
$ws.Cells.Item(111,14).Value = 1.5
$ws.Cells.Item(111,15).Value = 3.6
$ws.Cells.Item(111,16).Value = 5.75
$ws.Cells.Item(111,17).Value = -1
$ws.Cells.Item(111,18).Value = 1.9
$ws.Cells.Item(111,19).Value = 1.9
$ws.Cells.Item(111,20).Value = 2.25
$ws.Cells.Item(111,21).Value = 1.85
$ws.Cells.Item(111,22).Value = 1.95
$ws.Cells.Item(111,23).Value = 0.5
$ws.Cells.Item(111,24).Value = -1
$ws.Cells.Item(111,25).Value = -1
$ws.Cells.Item(111,26).Value = 0.8999999999999999
$ws.Cells.Item(111,27).Value = -1
$ws.Cells.Item(111,28).Value = -0.5
$ws.Cells.Item(111,29).Value = 0.475

# Row 112
$ws.Cells.Item(112,2).Value = 6865354
$ws.Cells.Item(112,3).Value = "Bosnia Herzegovina Premier Liga"
$ws.Cells.Item(112,4).Value = "Bosnia  Herzegovina Premier Liga"
$ws.Cells.Item(112,5).Value = 45339.375
$ws.Cells.Item(112,6).Value = "NK Igman Konjic"
$ws.Cells.Item(112,7).Value = "GOSK Gabela"
$ws.Cells.Item(112,8).Value = 1
$ws.Cells.Item(112,9).Value = 2
$ws.Cells.Item(112,10).Value = "A"
$ws.Cells.Item(112,11).Value = 1.8
$ws.Cells.Item(112,12).Value = 3.25
$ws.Cells.Item(112,13).Value = 4
$ws.Cells.Item(112,14).Value = 2.25
$ws.Cells.Item(112,15).Value = 3.1
$ws.Cells.Item(112,16).Value = 2.9
$ws.Cells.Item(112,17).Value = -0.25
$ws.Cells.Item(112,18).Value = 1.975
$ws.Cells.Item(112,19).Value = 1.825
$ws.Cells.Item(112,20).Value = 2.25
$ws.Cells.Item(112,21).Value = 1.875
$ws.Cells.Item(112,22).Value = 1.925
$ws.Cells.Item(112,23).Value = -1
$ws.Cells.Item(112,24).Value = -1
$ws.Cells.Item(112,25).Value = 1.9
$ws.Cells.Item(112,26).Value = -1
$ws.Cells.Item(112,27).Value = 0.825
$ws.Cells.Item(112,28).Value = 0.875
$ws.Cells.Item(112,29).Value = -1

# --- Append new row 156 (new fixture record) ---
# Copy formatting from the last existing row (155) so style indices (bold/border id col, date fmt) are reused
$ws.Cells.Item(155,1).Copy()
$ws.Cells.Item(156,1).PasteSpecial(-4122)
$ws.Cells.Item(155,5).Copy()
$ws.Cells.Item(156,5).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(156,1).Value = 154
$ws.Cells.Item(156,2).Value = 7952748
$ws.Cells.Item(156,3).Value = "Bosnia Herzegovina Premier Liga"
$ws.Cells.Item(156,4).Value = "Bosnia  Herzegovina Premier Liga"
$ws.Cells.Item(156,5).Value = 45396.33333333334
$ws.Cells.Item(156,6).Value = "NK Igman Konjic"
$ws.Cells.Item(156,7).Value = "NK Posusje"
$ws.Cells.Item(156,8).Value = 1
$ws.Cells.Item(156,9).Value = 1
$ws.Cells.Item(156,10).Value = "D"
$ws.Cells.Item(156,11).Value = 2.2
$ws.Cells.Item(156,12).Value = 3.4
$ws.Cells.Item(156,13).Value = 2.75
$ws.Cells.Item(156,14).Value = 2.05
$ws.Cells.Item(156,15).Value = 3.25
$ws.Cells.Item(156,16).Value = 3.1
$ws.Cells.Item(156,17).Value = -0.5
$ws.Cells.Item(156,18).Value = 2.1
$ws.Cells.Item(156,19).Value = 1.7
$ws.Cells.Item(156,20).Value = 2.25
$ws.Cells.Item(156,21).Value = 2
$ws.Cells.Item(156,22).Value = 1.8
$ws.Cells.Item(156,23).Value = -1
$ws.Cells.Item(156,24).Value = 2.25
$ws.Cells.Item(156,25).Value = -1
$ws.Cells.Item(156,26).Value = -1
$ws.Cells.Item(156,27).Value = 0.7
$ws.Cells.Item(156,28).Value = -0.5
$ws.Cells.Item(156,29).Value = 0.4
